# Integrated new AUT step: add a "takeScreenshot" keyword row to the
# Login sheet's keyword-driven test table, then leave that row/sheet
# selected (matches the author switching focus to the Login tab).

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login")

# Append the new row: Sr No = 5, Keywords = takeScreenshot
$wsLogin.Range("A6").Value = 5
$wsLogin.Range("B6").Value = "takeScreenshot"

# Grow the keyword table so the new row becomes part of it
$tbl = $wsLogin.ListObjects.Item(1)
$tbl.Resize($wsLogin.Range("A1:E6"))

# Make the Login sheet active with the new cell selected
$wsLogin.Activate()
$wsLogin.Range("B6").Select()

$wb.Save()
